$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 540.11
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "11.28%"

$ws.Range("C3").Value = 60.32
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "12.33%"

$ws.Range("C4").Value = 759.09
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "15.99%"

$ws.Range("C5").Value = 212.69
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "13.38%"

$ws.Range("C6").Value = 266.95
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "8.18%"

$ws.Range("C7").Value = 85.06
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "10.94%"
